$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'24.959.65"
$ws.Range("E2").Value = "  +2.04%  "

# Row 3
$ws.Range("D3").Value = "'1.700.11"
$ws.Range("E3").Value = "  +0.70%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'316.02"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("D7").Value = "'0.3980"
$ws.Range("E7").Value = "  +1.52%  "

# Row 8
$ws.Range("D8").Value = "'0.4025"
$ws.Range("E8").Value = "  -0.20%  "

# Row 9
$ws.Range("E9").Value = "  -1.34%  "

# Row 10
$ws.Range("D10").Value = "'53.26"
$ws.Range("E10").Value = "  +1.41%  "

# Row 11
$ws.Range("D11").Value = "'1.002"
$ws.Range("E11").Value = "  -0.07%  "

# Row 12
$ws.Range("D12").Value = "'0.08804"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13
$ws.Range("D13").Value = "'25.91"
$ws.Range("E13").Value = "  -2.13%  "

# Row 14
$ws.Range("D14").Value = "'7.467"
$ws.Range("E14").Value = "  +0.09%  "

# Row 15
$ws.Range("D15").Value = "'0.00001354"
$ws.Range("E15").Value = "  +0.74%  "

# Row 16
$ws.Range("D16").Value = "'7.960"
$ws.Range("E16").Value = "  -2.40%  "

# Row 17
$ws.Range("D17").Value = "'1.707.78"
$ws.Range("E17").Value = "  +1.30%  "

# Row 18
$ws.Range("D18").Value = "'95.59"
$ws.Range("E18").Value = "  -2.60%  "

# Row 19
$ws.Range("D19").Value = "'0.07190"
$ws.Range("E19").Value = "  -0.90%  "

# Row 20
$ws.Range("D20").Value = "'20.64"
$ws.Range("E20").Value = "  +2.10%  "

# Row 21
$ws.Range("D21").Value = "'7.317"
$ws.Range("E21").Value = "  +0.38%  "

# Row 22
$ws.Range("E22").Value = "  +0.06%  "

# Row 23
$ws.Range("D23").Value = "'14.35"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").Value = "'24.966.46"
$ws.Range("E24").Value = "  +2.12%  "

# Row 25
$ws.Range("D25").Value = "'2.383"
$ws.Range("E25").Value = "  +1.95%  "

# Row 26
$ws.Range("D26").Value = "'2.945"
$ws.Range("E26").Value = "  -2.80%  "

# Row 27
$ws.Range("E27").Value = "  +4.72%  "

# Row 28
$ws.Range("D28").Value = "'6.169"
$ws.Range("E28").Value = "  +15.24%  "

# Row 29
$ws.Range("D29").Value = "'162.60"
$ws.Range("E29").Value = "  -2.83%  "

# Row 30
$ws.Range("D30").Value = "'150.71"
$ws.Range("E30").Value = "  +8.98%  "

# Row 31
$ws.Range("D31").Value = "'8.355"
$ws.Range("E31").Value = "  -1.36%  "

# Row 32
$ws.Range("D32").Value = "'2.647"
$ws.Range("E32").Value = "  +26.97%  "

# Row 33
$ws.Range("D33").Value = "'1.895.46"
$ws.Range("E33").Value = "  +1.48%  "

# Row 34
$ws.Range("D34").Value = "'0.08541"
$ws.Range("E34").Value = "  -2.14%  "

# Row 35
$ws.Range("E35").Value = "  +4.99%  "

# Row 36
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'7.167"
$ws.Range("E36").Value = "  -1.04%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.038"
$ws.Range("E37").Value = "  -0.43%  "

# Row 38
$ws.Range("D38").Value = "'0.2879"
$ws.Range("E38").Value = "  +3.39%  "

# Row 39
$ws.Range("D39").Value = "'0.09601"
$ws.Range("E39").Value = "  +4.84%  "

# Row 40
$ws.Range("D40").Value = "'10.87"
$ws.Range("E40").Value = "  +0.70%  "

# Row 41
$ws.Range("D41").Value = "'0.8256"
$ws.Range("E41").Value = "  +2.20%  "

# Row 42
$ws.Range("D42").Value = "'14.01"
$ws.Range("E42").Value = "  -0.94%  "

# Row 43
$ws.Range("D43").Value = "'1.476"
$ws.Range("E43").Value = "  +0.42%  "

# Row 44
$ws.Range("D44").Value = "'17.25"
$ws.Range("E44").Value = "  -1.65%  "

# Row 45
$ws.Range("E45").Value = "  +1.18%  "

# Row 46
$ws.Range("D46").Value = "'0.7387"
$ws.Range("E46").Value = "  +1.81%  "

# Row 47
$ws.Range("E47").Value = "  -0.29%  "

# Row 48
$ws.Range("D48").Value = "'1.395"
$ws.Range("E48").Value = "  -1.07%  "

# Row 49
$ws.Range("D49").Value = "'0.08806"
$ws.Range("E49").Value = "  +8.53%  "

# Row 50
$ws.Range("E50").Value = "  +0.15%  "

# Row 51
$ws.Range("D51").Value = "'139.38"
$ws.Range("E51").Value = "  +0.04%  "
